$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold + "Text" (@) number format for the label cells in column A of the
# three SSB dataset rows (7, 9, 10) - matches cellXfs style index 2.
# (Multi-area Range() assignments only touch the first area in this host,
# so iterate .Areas explicitly.)
$boldTextCells = $ws.Range("A7,A9,A10")
foreach ($area in $boldTextCells.Areas) {
    $area.Font.Bold = $true
    $area.NumberFormat = "@"
}

# Bold (General number format) for the rest of the data in those same rows
# - matches cellXfs style index 3.
$boldGeneralCells = $ws.Range("B7:F7,B9:F9,B10:F10")
foreach ($area in $boldGeneralCells.Areas) {
    $area.Font.Bold = $true
}

# New row 17: a small "EIA" / "Oljepriser" dataset entry, using the same
# bold styles introduced above.
$ws.Range("A17").Value = "EIA"
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").NumberFormat = "@"

$ws.Range("B17").Value = "Oljepriser"
$ws.Range("B17").Font.Bold = $true

# Move the active selection to A8, matching the saved view state.
$ws.Range("A8").Select() | Out-Null
